$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.440.98"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.917.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4698"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2862"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06833"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.94"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +11.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.47"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07741"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.897.36"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.292"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6587"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "296.04"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.448.73"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007645"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.95"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.138.34"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.253"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.236"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.84"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.369"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.70"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.01%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.364"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.185"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.991"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05048"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7365"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.155"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02075"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.737"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.678"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.062"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "109.62"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8735"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.837"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4260"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "51.41"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +19.89%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.200"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.281"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2462"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +10.60%  "
